$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the previous row (row 4) down into the new row 5,
# so the new row matches the existing style pattern (s=1,0,2,2,0,0).
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Fill in the new review row's values.
$ws.Range("A5").Value = "com.singleton.strechy"
$ws.Range("B5").Value = "stretchy"
$ws.Range("C5").Value = "veredsnir12@gmail.com"
$ws.Range("D5").Value = "kevinkors122@gmail.com"
$ws.Range("E5").Value = "27/5/2019 15:60"
$ws.Range("F5").Value = "This is a tremendous playgame! A lot of cars to choose from. Free and offline game."

# Add the mailto hyperlinks for the email columns, mirroring the existing rows.
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")

# Adding the hyperlinks re-styles the cells with a "Hyperlink" look; restore
# the original formatting (style s=2) used by the other email cells.
$ws.Range("C4:D4").Copy()
$ws.Range("C5:D5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Match the author's final selection/active cell.
[void]$ws.Range("F5").Select()
